$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тесты")

# --- Row 42 (test #41): new "adaptive pool increase" experiment row ---
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = 100
$ws.Range("D42").Value = 25
$ws.Range("F42").Value = "Adaptive pool увеличен с (4,4) до (7,7)"
$ws.Range("G42").Value = "параметры теста 4"
$ws.Range("H42").Value = "Train IoU: 0.64, Val IoU: 0.64, Test IoU: 0.64."
$ws.Rows.Item(42).RowHeight = 30

# --- Row 43 (test #42): new "adaptive pool decrease" experiment row ---
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = 100
$ws.Range("D43").Value = 26
$ws.Range("F43").Value = "Adaptive pool уменьшен с (7,7) до (5,5)"
$ws.Range("G43").Value = "параметры теста 4"
$ws.Range("H43").Value = "Train IoU: 0.66, Val IoU: 0.67, Test IoU: 0.68."
$ws.Rows.Item(43).RowHeight = 30

# Move the active selection to H43 (matches the author's final cursor position)
$ws.Range("H43").Select()
